# Workbook / worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Row 1 (header): was a stray duplicate of the data row, fix it to the
# real column-name header used on the other property sheets, and extend it
# with the standard metadata columns (property_category .. index). ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Match the bold/bordered header formatting already used on B1:G1.
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

# --- Row 2 (data): B2:G2 keep their existing values; add the standard
# metadata columns so this sheet matches 土地/建物/股票. ---
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2013-12-27"
$ws.Range("K2").Value = "王惠美"
$ws.Range("L2").Value = 1729
$ws.Range("M2").Value = "tmp299c1"
$ws.Range("N2").Value = 30

# Match the plain/bordered data-row formatting already used on B2:G2.
$ws.Range("B2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)

$excel.CutCopyMode = $false
